# Weekly NYPD CompStat report refresh: new crime data collected.
# Updates the report header (volume/week-of text), and the Week-to-Date /
# 28-Day / Year-to-Date / 2-Year crime statistics table rows 15-28, 33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant used to copy a donor cell's number format/style
# onto a target cell without disturbing its value.
$xlPasteFormats = -4122

function Copy-CellFormat($srcRow, $srcCol, $dstRow, $dstCol) {
    $ws.Cells.Item($srcRow, $srcCol).Copy() | Out-Null
    $ws.Cells.Item($dstRow, $dstCol).PasteSpecial($xlPasteFormats) | Out-Null
}

# Set a cell to a numeric value while adopting the numeric style of a donor
# cell (used where the cell previously held a text placeholder, e.g. "0" /
# "***.*", and now needs to hold a real number).
function Set-NumberFromText($row, $col, $value, $donorRow, $donorCol) {
    Copy-CellFormat $donorRow $donorCol $row $col
    $ws.Cells.Item($row, $col).Value = $value
}

# Set a cell to a text value while adopting the text style of a donor cell
# (used where the cell previously held a real number and now needs to hold
# a text placeholder, e.g. "0" / "***.*").
function Set-TextFromNumber($row, $col, $text, $donorRow, $donorCol) {
    $ws.Cells.Item($row, $col).NumberFormat = "@"
    $ws.Cells.Item($row, $col).Value = $text
    Copy-CellFormat $donorRow $donorCol $row $col
}

# ---------------------------------------------------------------------
# Header: volume/number and report week text (rich-text shared strings).
# ---------------------------------------------------------------------

# A8 = "Volume 31   Number  32" -> "...33"
$ws.Cells.Item(8, 1).Characters(21, 2).Text = "33"

# C9 = "Report Covering the Week  8/5/2024  Through  8/11/2024"
#   -> "Report Covering the Week  8/12/2024  Through  8/18/2024"
# (replace the later substring first so the earlier offset stays valid)
$ws.Cells.Item(9, 3).Characters(46, 9).Text = "8/18/2024"
$ws.Cells.Item(9, 3).Characters(27, 8).Text = "8/12/2024"

# ---------------------------------------------------------------------
# Row 15 (Rape)
# ---------------------------------------------------------------------
$ws.Cells.Item(15, 12).Value = -30

# ---------------------------------------------------------------------
# Row 16 (Robbery)
# ---------------------------------------------------------------------
$ws.Cells.Item(16, 3).Value = 1
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 6).Value = 11
$ws.Cells.Item(16, 7).Value = 9
$ws.Cells.Item(16, 8).Value = 22.222222222222
$ws.Cells.Item(16, 9).Value = 65
$ws.Cells.Item(16, 10).Value = 71
$ws.Cells.Item(16, 11).Value = -8.450704225352
$ws.Cells.Item(16, 12).Value = -24.418604651162
$ws.Cells.Item(16, 13).Value = -24.418604651162
$ws.Cells.Item(16, 14).Value = -82.984293193717

# ---------------------------------------------------------------------
# Row 17 (Fel. Assault)
# ---------------------------------------------------------------------
$ws.Cells.Item(17, 3).Value = 3
$ws.Cells.Item(17, 4).Value = 3
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 6).Value = 10
$ws.Cells.Item(17, 7).Value = 17
$ws.Cells.Item(17, 8).Value = -41.176470588235
$ws.Cells.Item(17, 9).Value = 83
$ws.Cells.Item(17, 10).Value = 73
$ws.Cells.Item(17, 11).Value = 13.698630136986
$ws.Cells.Item(17, 12).Value = -15.306122448979
$ws.Cells.Item(17, 13).Value = 56.603773584905
$ws.Cells.Item(17, 14).Value = -50

# ---------------------------------------------------------------------
# Row 18 (Burglary)
# ---------------------------------------------------------------------
$ws.Cells.Item(18, 3).Value = 4
$ws.Cells.Item(18, 4).Value = 6
$ws.Cells.Item(18, 5).Value = -33.333333333333
$ws.Cells.Item(18, 6).Value = 10
$ws.Cells.Item(18, 7).Value = 19
$ws.Cells.Item(18, 8).Value = -47.368421052631
$ws.Cells.Item(18, 9).Value = 92
$ws.Cells.Item(18, 10).Value = 124
$ws.Cells.Item(18, 11).Value = -25.806451612903
$ws.Cells.Item(18, 12).Value = -40.645161290322
$ws.Cells.Item(18, 13).Value = -41.772151898734
$ws.Cells.Item(18, 14).Value = -87.239944521497

# ---------------------------------------------------------------------
# Row 19 (Gr. Larceny)
# ---------------------------------------------------------------------
$ws.Cells.Item(19, 3).Value = 10
$ws.Cells.Item(19, 4).Value = 17
$ws.Cells.Item(19, 5).Value = -41.176470588235
$ws.Cells.Item(19, 6).Value = 56
$ws.Cells.Item(19, 7).Value = 53
$ws.Cells.Item(19, 8).Value = 5.660377358490
$ws.Cells.Item(19, 9).Value = 403
$ws.Cells.Item(19, 10).Value = 428
$ws.Cells.Item(19, 11).Value = -5.841121495327
$ws.Cells.Item(19, 12).Value = 11.325966850828
$ws.Cells.Item(19, 13).Value = 137.058823529412
$ws.Cells.Item(19, 14).Value = 80.717488789237

# ---------------------------------------------------------------------
# Row 20 (G.L.A.)
# ---------------------------------------------------------------------
$ws.Cells.Item(20, 3).Value = 3
$ws.Cells.Item(20, 4).Value = 3
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = 9
$ws.Cells.Item(20, 7).Value = 20
$ws.Cells.Item(20, 8).Value = -55
$ws.Cells.Item(20, 9).Value = 61
$ws.Cells.Item(20, 10).Value = 107
$ws.Cells.Item(20, 11).Value = -42.990654205607
$ws.Cells.Item(20, 12).Value = -44.545454545454
$ws.Cells.Item(20, 13).Value = -35.789473684210
$ws.Cells.Item(20, 14).Value = -89.464594127806

# ---------------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------------
$ws.Cells.Item(21, 3).Value = 21
$ws.Cells.Item(21, 4).Value = 30
$ws.Cells.Item(21, 5).Value = -30
$ws.Cells.Item(21, 7).Value = 118
$ws.Cells.Item(21, 8).Value = -18.644067796610
$ws.Cells.Item(21, 9).Value = 711
$ws.Cells.Item(21, 10).Value = 806
$ws.Cells.Item(21, 11).Value = -11.786600496277
$ws.Cells.Item(21, 12).Value = -13.503649635036
$ws.Cells.Item(21, 13).Value = 26.287744227353
$ws.Cells.Item(21, 14).Value = -65.866538646183

# ---------------------------------------------------------------------
# Row 22 (Transit) - D/E/F/G/H swap between text placeholder and number
# ---------------------------------------------------------------------
Set-NumberFromText 22 4 1 22 9        # D22: "0" -> 1        (style like I22)
Set-NumberFromText 22 5 -100 22 11    # E22: "***.*" -> -100 (style like K22)
Set-TextFromNumber 22 6 "0" 22 3      # F22: 1 -> "0"        (style like C22)
Set-NumberFromText 22 7 1 22 9        # G22: "0" -> 1        (style like I22)
Set-NumberFromText 22 8 -100 22 11    # H22: "***.*" -> -100 (style like K22)
$ws.Cells.Item(22, 10).Value = 8
$ws.Cells.Item(22, 11).Value = 25
$ws.Cells.Item(22, 12).Value = 25

# ---------------------------------------------------------------------
# Row 23 (Housing) - D/E/F swap between text placeholder and number
# ---------------------------------------------------------------------
Set-TextFromNumber 23 4 "0" 23 3      # D23: 2 -> "0"      (style like C23)
Set-TextFromNumber 23 5 "***.*" 23 3  # E23: -100 -> "***.*" (style like C23)
Set-NumberFromText 23 6 2 23 9        # F23: "0" -> 2      (style like I23)
$ws.Cells.Item(23, 7).Value = 2
$ws.Cells.Item(23, 13).Value = -10

# ---------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------
$ws.Cells.Item(24, 3).Value = 22
$ws.Cells.Item(24, 4).Value = 17
$ws.Cells.Item(24, 5).Value = 29.411764705882
$ws.Cells.Item(24, 7).Value = 70
$ws.Cells.Item(24, 8).Value = 20
$ws.Cells.Item(24, 9).Value = 628
$ws.Cells.Item(24, 10).Value = 581
$ws.Cells.Item(24, 11).Value = 8.089500860585
$ws.Cells.Item(24, 12).Value = -1.412872841444
$ws.Cells.Item(24, 13).Value = 73.961218836565

# ---------------------------------------------------------------------
# Row 25 (Retail Theft)
# ---------------------------------------------------------------------
$ws.Cells.Item(25, 3).Value = 12
$ws.Cells.Item(25, 4).Value = 15
$ws.Cells.Item(25, 5).Value = -20
$ws.Cells.Item(25, 6).Value = 57
$ws.Cells.Item(25, 7).Value = 40
$ws.Cells.Item(25, 8).Value = 42.5
$ws.Cells.Item(25, 9).Value = 370
$ws.Cells.Item(25, 10).Value = 304
$ws.Cells.Item(25, 11).Value = 21.710526315789
$ws.Cells.Item(25, 12).Value = 0.817438692098

# ---------------------------------------------------------------------
# Row 26 (Misd. Assault)
# ---------------------------------------------------------------------
$ws.Cells.Item(26, 3).Value = 6
$ws.Cells.Item(26, 4).Value = 5
$ws.Cells.Item(26, 5).Value = 20
$ws.Cells.Item(26, 6).Value = 15
$ws.Cells.Item(26, 8).Value = -25
$ws.Cells.Item(26, 9).Value = 170
$ws.Cells.Item(26, 10).Value = 153
$ws.Cells.Item(26, 11).Value = 11.111111111111
$ws.Cells.Item(26, 12).Value = -3.954802259887
$ws.Cells.Item(26, 13).Value = 20.567375886524

# ---------------------------------------------------------------------
# Row 27 (UCR Rape*)
# ---------------------------------------------------------------------
$ws.Cells.Item(27, 12).Value = -33.333333333333

# ---------------------------------------------------------------------
# Row 28 (Other Sex Crimes) - C swaps from text placeholder to number
# ---------------------------------------------------------------------
Set-NumberFromText 28 3 2 28 4   # C28: "0" -> 2  (style like D28, a same-row number cell)
$ws.Cells.Item(28, 4).Value = 2
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 6
$ws.Cells.Item(28, 7).Value = 3
$ws.Cells.Item(28, 8).Value = 100
$ws.Cells.Item(28, 9).Value = 23
$ws.Cells.Item(28, 10).Value = 30
$ws.Cells.Item(28, 11).Value = -23.333333333333
$ws.Cells.Item(28, 12).Value = 21.052631578947

# ---------------------------------------------------------------------
# Row 33 (Traffic Fatalities) - G/H swap from number to text placeholder
# ---------------------------------------------------------------------
Set-TextFromNumber 33 7 "0" 33 3      # G33: 1 -> "0"        (style like C33)
Set-TextFromNumber 33 8 "***.*" 33 3  # H33: -100 -> "***.*" (style like C33)

$excel.CutCopyMode = 0
